# Ulazi i izlazi za Raspberry Pi Pico
# Consolidates split w:r runs (left over from incremental manual edits) into
# single runs with the same visible text, fixes the Pin(0, Pin.IN) row to
# read Pin(0, Pin.IN, Pin.PULL_DOWN), and drops the stray justification on
# that paragraph.

$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 1)
    if (-not $ok) {
        throw "Find failed for paragraph $paraIndex : '$findText'"
    }
}

# Row: Pin(16, Pin.OUT) ... (red 1)  -> merge the "(" "red" " 1)" runs
Replace-InParagraph 8 "(red 1)" "(red 1)"

# Row: Pin(0, Pin.IN) ... (button) -> becomes Pin(0, Pin.IN, Pin.PULL_DOWN)
$oldPin0Text = "Pin(0, Pin.IN)" + " " + "                                                  "
Replace-InParagraph 9 $oldPin0Text "Pin(0, Pin.IN, Pin.PULL_DOWN)     "
# Drop the <w:jc w:val="both"/> on that paragraph's pPr
$d.Paragraphs.Item(9).Format.Alignment = 0

# Row: Pin(17, Pin.OUT) ... (red 2)
Replace-InParagraph 11 "Pin(17, Pin.OUT)" "Pin(17, Pin.OUT)"
Replace-InParagraph 11 "(red 2)" "(red 2)"

# Row: Pin(20, Pin.IN, Pin.PULL_DOWN) ... (kol 1)
Replace-InParagraph 12 "(kol 1)" "(kol 1)"

# Row: Pin(18, Pin.OUT) ... (red 3)
Replace-InParagraph 14 "Pin(18, Pin.OUT)" "Pin(18, Pin.OUT)"
Replace-InParagraph 14 "(red 3)" "(red 3)"

# Row: Pin(21, Pin.IN, Pin.PULL_DOWN) ... (kol 2)
Replace-InParagraph 15 "Pin(21, Pin.IN, Pin.PULL_DOWN)" "Pin(21, Pin.IN, Pin.PULL_DOWN)"
Replace-InParagraph 15 "(kol 2)" "(kol 2)"

# Row: Pin(19, Pin.OUT) ... (red 4)
Replace-InParagraph 17 "Pin(19, Pin.OUT)" "Pin(19, Pin.OUT)"
Replace-InParagraph 17 "(red 4)" "(red 4)"

# Row: Pin(22, Pin.IN, Pin.PULL_DOWN) ... (kol 3)
Replace-InParagraph 18 "Pin(22, Pin.IN, Pin.PULL_DOWN)" "Pin(22, Pin.IN, Pin.PULL_DOWN)"
Replace-InParagraph 18 "(kol 3)" "(kol 3)"

# Row: Pin(7, Pin.OUT) ... (crvena LED)
Replace-InParagraph 20 "(crvena LED)" "(crvena LED)"

# Row: Pin(26, Pin.IN, Pin.PULL_DOWN) ... (kol 4)
Replace-InParagraph 21 "Pin(26, Pin.IN, Pin.PULL_DOWN)" "Pin(26, Pin.IN, Pin.PULL_DOWN)"
Replace-InParagraph 21 "(kol 4)" "(kol 4)"

# Row: Pin(6, Pin.OUT) ... (zelena LED)
Replace-InParagraph 23 "Pin(6, Pin.OUT)" "Pin(6, Pin.OUT)"
Replace-InParagraph 23 "(zelena LED)" "(zelena LED)"

# I2C table header cells: "I2C" + " izlazi" / "I2C" + " ulazi" -> single runs
Replace-InParagraph 55 "I2C izlazi" "I2C izlazi"
Replace-InParagraph 56 "I2C ulazi" "I2C ulazi"

Write-Output "done"
